$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column L, matching header style of K1 (bold/centered)
$ws.Range("L1").Value = "break_on_off"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# Fill L2:L73 with 0 by default
$ws.Range("L2:L73").Value = 0

# Set specific rows to 1 (break on/off flags)
$ws.Range("L19").Value = 1
$ws.Range("L37").Value = 1
$ws.Range("L54").Value = 1

# Update selection to reflect the new active cell from the diff
$ws.Range("O13").Select() | Out-Null
